# Adds a new "2022-Q4" sheet (fund holding detail) right after "总计",
# and updates the "总计" summary sheet with a new row for 2022-Q4,
# shifting the existing quarters down.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$val) {
    # Force a value to be stored as text (preserves leading zeros / trailing
    # zeros in numeric-looking strings), then strip the implicit "text"
    # number format Excel applies so the cell keeps the default style.
    $range.Value = "'" + $val
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet by duplicating "2022-Q3" (so it
#    inherits identical column layout / header style / borders), then
#    overwrite its data with the 2022-Q4 figures.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $summary)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The source sheet only had 2 data rows; 2022-Q4 needs 3, so extend the
# table by copying row 3's formatting down into row 4.
$q4.Range("A3:H3").Copy()
$q4.Range("A4:H4").PasteSpecial(-4122)

# Row 2
$q4.Range("A2").Value = 0
Set-TextValue $q4.Range("B2") "168501"
$q4.Range("C2").Value = "北信瑞丰产业升级多策略混合"
Set-TextValue $q4.Range("D2") "1.47"
Set-TextValue $q4.Range("E2") "93.49"
Set-TextValue $q4.Range("F2") "8.03"
Set-TextValue $q4.Range("G2") "0.1180"
$q4.Range("H2").Value = 1

# Row 3
$q4.Range("A3").Value = 1
Set-TextValue $q4.Range("B3") "009954"
$q4.Range("C3").Value = "北信瑞丰优选成长股票"
Set-TextValue $q4.Range("D3") "0.54"
Set-TextValue $q4.Range("E3") "93.50"
Set-TextValue $q4.Range("F3") "7.98"
Set-TextValue $q4.Range("G3") "0.0431"
$q4.Range("H3").Value = 4

# Row 4
$q4.Range("A4").Value = 2
Set-TextValue $q4.Range("B4") "001829"
$q4.Range("C4").Value = "北信瑞丰中国智造主题灵活配置混合"
Set-TextValue $q4.Range("D4") "0.28"
Set-TextValue $q4.Range("E4") "93.54"
Set-TextValue $q4.Range("F4") "7.96"
Set-TextValue $q4.Range("G4") "0.0223"
$q4.Range("H4").Value = 3

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a 2022-Q4 row at the top of
#    the data (row 2) and push the rest down, re-numbering the index
#    column (A) and appending the row that falls off the bottom (2020-Q4
#    moves from row 9 to row 10).
# ---------------------------------------------------------------------
$dates  = @("2022-Q4","2022-Q3","2022-Q2","2022-Q1","2021-Q4","2021-Q3","2021-Q2","2021-Q1","2020-Q4")
$counts = @(3,2,6,6,3,6,4,2,1)
$values = @(0.18,0.04,0.29,0.11,0.13,0.14,0.09,0.02,0.02)

# Copy formatting of the last existing row (row 9) down into the new
# row 10 before filling it with data.
$summary.Range("A9:D9").Copy()
$summary.Range("A10:D10").PasteSpecial(-4122)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $summary.Cells.Item($row, 1).Value = $i
    $summary.Cells.Item($row, 2).Value = $dates[$i]
    $summary.Cells.Item($row, 3).Value = $counts[$i]
    $summary.Cells.Item($row, 4).Value = $values[$i]
}

# ---------------------------------------------------------------------
# 3. Leave the workbook focused on "总计", matching the original view.
# ---------------------------------------------------------------------
$summary.Activate()
